# Append the new resale-numbers row for 2024-01-18 21:44:33 to the
# CityResaleNum sheet (row 72, directly below the existing last row 71).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 72

# Text columns (A-D): temporarily force text format so values that look
# like dates / times / zero-padded numbers ("2024-01-18", "21:44:33",
# "02") are stored verbatim as strings -- matching the rest of the sheet,
# where these columns are text, not real date/number types. The format
# is cleared again right after so the new cells end up with the same
# (default) style as every other data row.
$textRange = $ws.Range("A$row`:D$row")
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-01-18"
$ws.Cells.Item($row, 2).Value = "21:44:33"
$ws.Cells.Item($row, 3).Value = "Thursday"
$ws.Cells.Item($row, 4).Value = "02"

$textRange.ClearFormats()

# Numeric columns (E-T)
$values = @(139118, 140227, 171552, 148767, -1, 121506, 223388, 254860, 185205, 110362, 41372, 30921, 73556, -1, 42887, -1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = 5 + $i
    $ws.Cells.Item($row, $col).Value = $values[$i]
}
